# Auto-generated edit script reproducing the committed diff.
# Strategy: rewrite every string-valued cell in the precise operation order required
# so the regenerated shared-string table matches the target order exactly, since the
# runtime (like Excel) appends newly-introduced unique strings to the table in the order
# that Range.Value assignments occur (not sheet position).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: column B (HKL label) rows 2-31, in order -> builds shared strings 0..29
$ws.Cells.Item(2,2).Value = 'HKL'
$ws.Cells.Item(3,2).Value = 'Spiral5'
$ws.Cells.Item(4,2).Value = 'Holden'
$ws.Cells.Item(5,2).Value = 'Rizzie Spiral'
$ws.Cells.Item(6,2).Value = 'RotRing OmegaMax-90'
$ws.Cells.Item(7,2).Value = 'Equal Angle'
$ws.Cells.Item(8,2).Value = 'Tilt Rotate'
$ws.Cells.Item(9,2).Value = 'CLR'
$ws.Cells.Item(10,2).Value = 'Rizzie Hex'
$ws.Cells.Item(11,2).Value = 'Matthies Hex'
$ws.Cells.Item(12,2).Value = 'Tilt Rotate_Partial'
$ws.Cells.Item(13,2).Value = 'RotRing OmegaMax-60'
$ws.Cells.Item(14,2).Value = 'Equal Angle_Partial'
$ws.Cells.Item(15,2).Value = 'Rizzie Hex_Partial'
$ws.Cells.Item(16,2).Value = 'ND Single'
$ws.Cells.Item(17,2).Value = 'RD Single'
$ws.Cells.Item(18,2).Value = 'TD Single'
$ws.Cells.Item(19,2).Value = 'Morris Single'
$ws.Cells.Item(20,2).Value = 'Ring Perpendicular to ND'
$ws.Cells.Item(21,2).Value = 'Ring Perpendicular to RD'
$ws.Cells.Item(22,2).Value = 'Ring Perpendicular to TD'
$ws.Cells.Item(23,2).Value = 'OffsetFTD'
$ws.Cells.Item(24,2).Value = 'OffsetATD'
$ws.Cells.Item(25,2).Value = 'OffsetF45'
$ws.Cells.Item(26,2).Value = 'OffsetA45'
$ws.Cells.Item(27,2).Value = 'OffsetFRD'
$ws.Cells.Item(28,2).Value = 'OffsetARD'
$ws.Cells.Item(29,2).Value = 'Gaussian Quadrature'
$ws.Cells.Item(30,2).Value = 'Michael-CCHex'
$ws.Cells.Item(31,2).Value = 'Michael-SNHex'

# --- Step 2: row 2 (bracket labels) columns C-W, in order -> builds shared strings 30..50
$ws.Cells.Item(2,3).Value = '[5, 1, 1]'
$ws.Cells.Item(2,4).Value = '[4, 2, 2]'
$ws.Cells.Item(2,5).Value = '[3, 1, 1]'
$ws.Cells.Item(2,6).Value = '[3, 3, 1]'
$ws.Cells.Item(2,7).Value = '[2, 2, 2]'
$ws.Cells.Item(2,8).Value = '[1, 1, 1]'
$ws.Cells.Item(2,9).Value = '[3, 3, 3]'
$ws.Cells.Item(2,10).Value = '[2, 2, 0]'
$ws.Cells.Item(2,11).Value = '[2, 0, 0]'
$ws.Cells.Item(2,12).Value = '[4, 0, 0]'
$ws.Cells.Item(2,13).Value = '[4, 2, 0]'
$ws.Cells.Item(2,14).Value = '1Pair-A'
$ws.Cells.Item(2,15).Value = '1Pair-B'
$ws.Cells.Item(2,16).Value = '2Pairs-A'
$ws.Cells.Item(2,17).Value = '2Pairs-B'
$ws.Cells.Item(2,18).Value = '3Pairs-A'
$ws.Cells.Item(2,19).Value = '3Pairs-B'
$ws.Cells.Item(2,20).Value = '3Pairs-C'
$ws.Cells.Item(2,21).Value = '4Pairs'
$ws.Cells.Item(2,22).Value = '5A4F'
$ws.Cells.Item(2,23).Value = 'MaxUnique'

# --- Step 3: column A (numeric index) rows 2-31
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(12,1).Value = 10
$ws.Cells.Item(13,1).Value = 11
$ws.Cells.Item(14,1).Value = 12
$ws.Cells.Item(15,1).Value = 13
$ws.Cells.Item(16,1).Value = 14
$ws.Cells.Item(17,1).Value = 15
$ws.Cells.Item(18,1).Value = 16
$ws.Cells.Item(19,1).Value = 17
$ws.Cells.Item(20,1).Value = 18
$ws.Cells.Item(21,1).Value = 19
$ws.Cells.Item(22,1).Value = 20
$ws.Cells.Item(23,1).Value = 21
$ws.Cells.Item(24,1).Value = 22
$ws.Cells.Item(25,1).Value = 23
$ws.Cells.Item(26,1).Value = 24
$ws.Cells.Item(27,1).Value = 25
$ws.Cells.Item(28,1).Value = 26
$ws.Cells.Item(29,1).Value = 27
$ws.Cells.Item(30,1).Value = 28
$ws.Cells.Item(31,1).Value = 29

# --- Step 4: row 1 numeric header (unchanged by the edit, rewritten for completeness)
$ws.Cells.Item(1,2).Value = 0
$ws.Cells.Item(1,3).Value = 1
$ws.Cells.Item(1,4).Value = 2
$ws.Cells.Item(1,5).Value = 3
$ws.Cells.Item(1,6).Value = 4
$ws.Cells.Item(1,7).Value = 5
$ws.Cells.Item(1,8).Value = 6
$ws.Cells.Item(1,9).Value = 7
$ws.Cells.Item(1,10).Value = 8
$ws.Cells.Item(1,11).Value = 9
$ws.Cells.Item(1,12).Value = 10
$ws.Cells.Item(1,13).Value = 11
$ws.Cells.Item(1,14).Value = 12
$ws.Cells.Item(1,15).Value = 13
$ws.Cells.Item(1,16).Value = 14
$ws.Cells.Item(1,17).Value = 15
$ws.Cells.Item(1,18).Value = 16
$ws.Cells.Item(1,19).Value = 17
$ws.Cells.Item(1,20).Value = 18
$ws.Cells.Item(1,21).Value = 19
$ws.Cells.Item(1,22).Value = 20
$ws.Cells.Item(1,23).Value = 21

# --- Step 5: row 3 data (unchanged by the edit, rewritten for completeness)
$ws.Cells.Item(3,3).Value = 0.9984658230784822
$ws.Cells.Item(3,4).Value = 1.000551785957545
$ws.Cells.Item(3,5).Value = 0.9995037810448227
$ws.Cells.Item(3,6).Value = 1.000911040086609
$ws.Cells.Item(3,7).Value = 1.001336596072072
$ws.Cells.Item(3,8).Value = 1.001336596072072
$ws.Cells.Item(3,9).Value = 1.001336596072072
$ws.Cells.Item(3,10).Value = 1.000723163314588
$ws.Cells.Item(3,11).Value = 0.9976431787010642
$ws.Cells.Item(3,12).Value = 0.9976431787010642
$ws.Cells.Item(3,13).Value = 0.9995736563911292
$ws.Cells.Item(3,14).Value = 1.001336596072072
$ws.Cells.Item(3,15).Value = 1.000723163314588
$ws.Cells.Item(3,16).Value = 0.9991831710078262
$ws.Cells.Item(3,17).Value = 1.000113472179705
$ws.Cells.Item(3,18).Value = 0.9999009793625747
$ws.Cells.Item(3,19).Value = 0.9992900410201583
$ws.Cells.Item(3,20).Value = 0.9999009793625747
$ws.Cells.Item(3,21).Value = 0.9998016797831367
$ws.Cells.Item(3,22).Value = 1.000108663040924
$ws.Cells.Item(3,23).Value = 0.9998386280807889

# --- Step 6: rows 4-31 data columns C:W (rows 4,5 are new simulation results;
#             rows 6-31 reuse the data previously on rows 4-29, shifted down by 2)
$ws.Cells.Item(4,3).Value = 0.9441377423279577
$ws.Cells.Item(4,4).Value = 1.014305408482217
$ws.Cells.Item(4,5).Value = 0.9816364538187327
$ws.Cells.Item(4,6).Value = 1.035680405775157
$ws.Cells.Item(4,7).Value = 1.031065157713039
$ws.Cells.Item(4,8).Value = 1.031065157713039
$ws.Cells.Item(4,9).Value = 1.031065157713039
$ws.Cells.Item(4,10).Value = 1.038593802438887
$ws.Cells.Item(4,11).Value = 0.9114039102761753
$ws.Cells.Item(4,12).Value = 0.9114039102761753
$ws.Cells.Item(4,13).Value = 0.9909873834273307
$ws.Cells.Item(4,14).Value = 1.031065157713039
$ws.Cells.Item(4,15).Value = 1.038593802438887
$ws.Cells.Item(4,16).Value = 0.974998856357531
$ws.Cells.Item(4,17).Value = 1.01011512812881
$ws.Cells.Item(4,18).Value = 0.9936876234760336
$ws.Cells.Item(4,19).Value = 0.9772113888445982
$ws.Cells.Item(4,20).Value = 0.9936876234760336
$ws.Cells.Item(4,21).Value = 0.9906748310617084
$ws.Cells.Item(4,22).Value = 0.9987528963919743
$ws.Cells.Item(4,23).Value = 0.9934762830324368
$ws.Cells.Item(5,3).Value = 0.7639092835037746
$ws.Cells.Item(5,4).Value = 1.075482861433785
$ws.Cells.Item(5,5).Value = 0.9037454846467128
$ws.Cells.Item(5,6).Value = 1.148241365742549
$ws.Cells.Item(5,7).Value = 1.239354499233444
$ws.Cells.Item(5,8).Value = 1.239354499233444
$ws.Cells.Item(5,9).Value = 1.239354499233444
$ws.Cells.Item(5,10).Value = 1.117048426297463
$ws.Cells.Item(5,11).Value = 0.6697065804066629
$ws.Cells.Item(5,12).Value = 0.6697065804066629
$ws.Cells.Item(5,13).Value = 0.9220546002350487
$ws.Cells.Item(5,14).Value = 1.239354499233444
$ws.Cells.Item(5,15).Value = 1.117048426297463
$ws.Cells.Item(5,16).Value = 0.8933775033520628
$ws.Cells.Item(5,17).Value = 1.010396955472088
$ws.Cells.Item(5,18).Value = 1.008703168645856
$ws.Cells.Item(5,19).Value = 0.8968334971169462
$ws.Cells.Item(5,20).Value = 1.008703168645856
$ws.Cells.Item(5,21).Value = 0.9824637476460705
$ws.Cells.Item(5,22).Value = 1.033841897963545
$ws.Cells.Item(5,23).Value = 0.97994288768743
$ws.Cells.Item(6,3).Value = 0.9096489231674431
$ws.Cells.Item(6,4).Value = 1.031882838794195
$ws.Cells.Item(6,5).Value = 0.9658390619108757
$ws.Cells.Item(6,6).Value = 1.053698775839675
$ws.Cells.Item(6,7).Value = 1.096617829388203
$ws.Cells.Item(6,8).Value = 1.096617829388203
$ws.Cells.Item(6,9).Value = 1.096617829388203
$ws.Cells.Item(6,10).Value = 1.038181789542507
$ws.Cells.Item(6,11).Value = 0.8681785816471163
$ws.Cells.Item(6,12).Value = 0.8681785816471163
$ws.Cells.Item(6,13).Value = 0.9704943620758754
$ws.Cells.Item(6,14).Value = 1.096617829388203
$ws.Cells.Item(6,15).Value = 1.038181789542507
$ws.Cells.Item(6,16).Value = 0.9531801855948117
$ws.Cells.Item(6,17).Value = 1.002010425726691
$ws.Cells.Item(6,18).Value = 1.000992733525942
$ws.Cells.Item(6,19).Value = 0.9573998110334996
$ws.Cells.Item(6,20).Value = 1.000992733525942
$ws.Cells.Item(6,21).Value = 0.9922043156221754
$ws.Cells.Item(6,22).Value = 1.013087018375381
$ws.Cells.Item(6,23).Value = 0.9918177702957363
$ws.Cells.Item(7,3).Value = 0.9291179567435147
$ws.Cells.Item(7,4).Value = 1.02012026463977
$ws.Cells.Item(7,5).Value = 0.9765108020533143
$ws.Cells.Item(7,6).Value = 1.044455468609508
$ws.Cells.Item(7,7).Value = 1.046519056181554
$ws.Cells.Item(7,8).Value = 1.046519056181554
$ws.Cells.Item(7,9).Value = 1.046519056181554
$ws.Cells.Item(7,10).Value = 1.044585316037463
$ws.Cells.Item(7,11).Value = 0.889005780237754
$ws.Cells.Item(7,12).Value = 0.889005780237754
$ws.Cells.Item(7,13).Value = 0.9859885573487034
$ws.Cells.Item(7,14).Value = 1.046519056181554
$ws.Cells.Item(7,15).Value = 1.044585316037463
$ws.Cells.Item(7,16).Value = 0.9667955481376083
$ws.Cells.Item(7,17).Value = 1.010548059045388
$ws.Cells.Item(7,18).Value = 0.9933700508189237
$ws.Cells.Item(7,19).Value = 0.9700339661095102
$ws.Cells.Item(7,20).Value = 0.9933700508189237
$ws.Cells.Item(7,21).Value = 0.9891552386275213
$ws.Cells.Item(7,22).Value = 1.000628002138328
$ws.Cells.Item(7,23).Value = 0.9920379002314477
$ws.Cells.Item(8,3).Value = 0.7634006996762543
$ws.Cells.Item(8,4).Value = 1.057842211715283
$ws.Cells.Item(8,5).Value = 0.9234264298596055
$ws.Cells.Item(8,6).Value = 1.155071307564205
$ws.Cells.Item(8,7).Value = 1.109240307005956
$ws.Cells.Item(8,8).Value = 1.109240307005956
$ws.Cells.Item(8,9).Value = 1.109240307005956
$ws.Cells.Item(8,10).Value = 1.176182004654559
$ws.Cells.Item(8,11).Value = 0.6262722612710234
$ws.Cells.Item(8,12).Value = 0.6262722612710234
$ws.Cells.Item(8,13).Value = 0.9644777645922072
$ws.Cells.Item(8,14).Value = 1.109240307005956
$ws.Cells.Item(8,15).Value = 1.176182004654559
$ws.Cells.Item(8,16).Value = 0.9012271329627914
$ws.Cells.Item(8,17).Value = 1.049804217257083
$ws.Cells.Item(8,18).Value = 0.9705648576438461
$ws.Cells.Item(8,19).Value = 0.9086268985950628
$ws.Cells.Item(8,20).Value = 0.9705648576438461
$ws.Cells.Item(8,21).Value = 0.9587802506977861
$ws.Cells.Item(8,22).Value = 0.9888722619594199
$ws.Cells.Item(8,23).Value = 0.9719891232923867
$ws.Cells.Item(9,3).Value = 0.9933863362670049
$ws.Cells.Item(9,4).Value = 1.003657232508415
$ws.Cells.Item(9,5).Value = 0.9976823095484235
$ws.Cells.Item(9,6).Value = 1.003210274467359
$ws.Cells.Item(9,7).Value = 1.0112195200874
$ws.Cells.Item(9,8).Value = 1.0112195200874
$ws.Cells.Item(9,9).Value = 1.0112195200874
$ws.Cells.Item(9,10).Value = 0.9998611476983256
$ws.Cells.Item(9,11).Value = 0.9905126237497666
$ws.Cells.Item(9,12).Value = 0.9905126237497666
$ws.Cells.Item(9,13).Value = 0.9965419623808349
$ws.Cells.Item(9,14).Value = 1.0112195200874
$ws.Cells.Item(9,15).Value = 0.9998611476983256
$ws.Cells.Item(9,16).Value = 0.995186885724046
$ws.Cells.Item(9,17).Value = 0.9987717286233746
$ws.Cells.Item(9,18).Value = 1.000531097178498
$ws.Cells.Item(9,19).Value = 0.9960186936655052
$ws.Cells.Item(9,20).Value = 1.000531097178498
$ws.Cells.Item(9,21).Value = 0.999818900270979
$ws.Cells.Item(9,22).Value = 1.002099024234263
$ws.Cells.Item(9,23).Value = 0.9995089258384413
$ws.Cells.Item(10,3).Value = 0.9994071531897171
$ws.Cells.Item(10,4).Value = 1.000282921522341
$ws.Cells.Item(10,5).Value = 1.000016789505508
$ws.Cells.Item(10,6).Value = 1.000272979213183
$ws.Cells.Item(10,7).Value = 1.000135937594402
$ws.Cells.Item(10,8).Value = 1.000135937594402
$ws.Cells.Item(10,9).Value = 1.000135937594402
$ws.Cells.Item(10,10).Value = 1.000278179579533
$ws.Cells.Item(10,11).Value = 0.9987035572099884
$ws.Cells.Item(10,12).Value = 0.9987035572099884
$ws.Cells.Item(10,13).Value = 0.9999921720315632
$ws.Cells.Item(10,14).Value = 1.000135937594402
$ws.Cells.Item(10,15).Value = 1.000278179579533
$ws.Cells.Item(10,16).Value = 0.9994908683947605
$ws.Cells.Item(10,17).Value = 1.00014748454252
$ws.Cells.Item(10,18).Value = 0.9997058914613076
$ws.Cells.Item(10,19).Value = 0.9996661754316764
$ws.Cells.Item(10,20).Value = 0.9997058914613076
$ws.Cells.Item(10,21).Value = 0.9997836159723577
$ws.Cells.Item(10,22).Value = 0.9998540802967664
$ws.Cells.Item(10,23).Value = 0.9998862112307794
$ws.Cells.Item(11,3).Value = 0.9892523666078221
$ws.Cells.Item(11,4).Value = 1.005502931367238
$ws.Cells.Item(11,5).Value = 0.9962814640006342
$ws.Cells.Item(11,6).Value = 1.005348167382216
$ws.Cells.Item(11,7).Value = 1.016764531306731
$ws.Cells.Item(11,8).Value = 1.016764531306731
$ws.Cells.Item(11,9).Value = 1.016764531306731
$ws.Cells.Item(11,10).Value = 1.000729919673267
$ws.Cells.Item(11,11).Value = 0.984219402780832
$ws.Cells.Item(11,12).Value = 0.984219402780832
$ws.Cells.Item(11,13).Value = 0.9949813394367369
$ws.Cells.Item(11,14).Value = 1.016764531306731
$ws.Cells.Item(11,15).Value = 1.000729919673267
$ws.Cells.Item(11,16).Value = 0.9924746612270493
$ws.Cells.Item(11,17).Value = 0.9985056918369504
$ws.Cells.Item(11,18).Value = 1.000571284586943
$ws.Cells.Item(11,19).Value = 0.993743595484911
$ws.Cells.Item(11,20).Value = 1.000571284586943
$ws.Cells.Item(11,21).Value = 0.9994988294403658
$ws.Cells.Item(11,22).Value = 1.002951969813639
$ws.Cells.Item(11,23).Value = 0.9991350153194345
$ws.Cells.Item(12,3).Value = 0.7595059371072382
$ws.Cells.Item(12,4).Value = 1.058158943493794
$ws.Cells.Item(12,5).Value = 0.9220707749319768
$ws.Cells.Item(12,6).Value = 1.157901533088965
$ws.Cells.Item(12,7).Value = 1.109284061402891
$ws.Cells.Item(12,8).Value = 1.109284061402891
$ws.Cells.Item(12,9).Value = 1.109284061402891
$ws.Cells.Item(12,10).Value = 1.18039527554535
$ws.Cells.Item(12,11).Value = 0.6199667690915767
$ws.Cells.Item(12,12).Value = 0.6199667690915767
$ws.Cells.Item(12,13).Value = 0.9645427618426494
$ws.Cells.Item(12,14).Value = 1.109284061402891
$ws.Cells.Item(12,15).Value = 1.18039527554535
$ws.Cells.Item(12,16).Value = 0.9001810223184632
$ws.Cells.Item(12,17).Value = 1.051233025238663
$ws.Cells.Item(12,18).Value = 0.9698820353466059
$ws.Cells.Item(12,19).Value = 0.9074776065229676
$ws.Cells.Item(12,20).Value = 0.9698820353466059
$ws.Cells.Item(12,21).Value = 0.9579292202429486
$ws.Cells.Item(12,22).Value = 0.9882001884749372
$ws.Cells.Item(12,23).Value = 0.9714782570630551
$ws.Cells.Item(13,3).Value = 0.9575432156080381
$ws.Cells.Item(13,4).Value = 1.017155768665251
$ws.Cells.Item(13,5).Value = 0.9809869377061212
$ws.Cells.Item(13,6).Value = 1.025359805803733
$ws.Cells.Item(13,7).Value = 1.060459748977564
$ws.Cells.Item(13,8).Value = 1.060459748977564
$ws.Cells.Item(13,9).Value = 1.060459748977564
$ws.Cells.Item(13,10).Value = 1.011941675500995
$ws.Cells.Item(13,11).Value = 0.945731409254036
$ws.Cells.Item(13,12).Value = 0.945731409254036
$ws.Cells.Item(13,13).Value = 0.9797661661113223
$ws.Cells.Item(13,14).Value = 1.060459748977564
$ws.Cells.Item(13,15).Value = 1.011941675500995
$ws.Cells.Item(13,16).Value = 0.9788365423775156
$ws.Cells.Item(13,17).Value = 0.9964643066035581
$ws.Cells.Item(13,18).Value = 1.006044277910865
$ws.Cells.Item(13,19).Value = 0.9795533408203841
$ws.Cells.Item(13,20).Value = 1.006044277910865
$ws.Cells.Item(13,21).Value = 0.9997799428596792
$ws.Cells.Item(13,22).Value = 1.011915904083256
$ws.Cells.Item(13,23).Value = 0.9973680909533826
$ws.Cells.Item(14,3).Value = 0.9266773381263157
$ws.Cells.Item(14,4).Value = 1.016713702557894
$ws.Cells.Item(14,5).Value = 0.9756137258631566
$ws.Cells.Item(14,6).Value = 1.047966981273682
$ws.Cells.Item(14,7).Value = 1.034264758463158
$ws.Cells.Item(14,8).Value = 1.034264758463158
$ws.Cells.Item(14,9).Value = 1.034264758463158
$ws.Cells.Item(14,10).Value = 1.055486558589472
$ws.Cells.Item(14,11).Value = 0.8836260383368416
$ws.Cells.Item(14,12).Value = 0.8836260383368416
$ws.Cells.Item(14,13).Value = 0.9901207649684197
$ws.Cells.Item(14,14).Value = 1.034264758463158
$ws.Cells.Item(14,15).Value = 1.055486558589472
$ws.Cells.Item(14,16).Value = 0.9695562984631567
$ws.Cells.Item(14,17).Value = 1.015550142226314
$ws.Cells.Item(14,18).Value = 0.9911257851298236
$ws.Cells.Item(14,19).Value = 0.9715754409298234
$ws.Cells.Item(14,20).Value = 0.9911257851298236
$ws.Cells.Item(14,21).Value = 0.9872477703131569
$ws.Cells.Item(14,22).Value = 0.9966511679431571
$ws.Cells.Item(14,23).Value = 0.9913087335223674
$ws.Cells.Item(15,3).Value = 1.043487418475062
$ws.Cells.Item(15,4).Value = 0.9848227026736269
$ws.Cells.Item(15,5).Value = 1.014764920246325
$ws.Cells.Item(15,6).Value = 0.9740807280127584
$ws.Cells.Item(15,7).Value = 0.9603851835183232
$ws.Cells.Item(15,8).Value = 0.9603851835183232
$ws.Cells.Item(15,9).Value = 0.9603851835183232
$ws.Cells.Item(15,10).Value = 0.9795660473333295
$ws.Cells.Item(15,11).Value = 1.066328163559763
$ws.Cells.Item(15,12).Value = 1.066328163559763
$ws.Cells.Item(15,13).Value = 1.012215343670952
$ws.Cells.Item(15,14).Value = 0.9603851835183232
$ws.Cells.Item(15,15).Value = 0.9795660473333295
$ws.Cells.Item(15,16).Value = 1.022947105446546
$ws.Cells.Item(15,17).Value = 0.9971654837898274
$ws.Cells.Item(15,18).Value = 1.002093131470472
$ws.Cells.Item(15,19).Value = 1.020219710379806
$ws.Cells.Item(15,20).Value = 1.002093131470472
$ws.Cells.Item(15,21).Value = 1.005261078664435
$ws.Cells.Item(15,22).Value = 0.9962858996352129
$ws.Cells.Item(15,23).Value = 1.004456313436267
$ws.Cells.Item(16,3).Value = 0.5835660400000007
$ws.Cells.Item(16,4).Value = 1.100097199999999
$ws.Cells.Item(16,5).Value = 0.8654664700000008
$ws.Cells.Item(16,6).Value = 1.274408899999999
$ws.Cells.Item(16,7).Value = 1.183405100000001
$ws.Cells.Item(16,8).Value = 1.183405100000001
$ws.Cells.Item(16,9).Value = 1.183405100000001
$ws.Cells.Item(16,10).Value = 1.315509199999998
$ws.Cells.Item(16,11).Value = 0.3422105699999997
$ws.Cells.Item(16,12).Value = 0.3422105699999997
$ws.Cells.Item(16,13).Value = 0.9392734700000017
$ws.Cells.Item(16,14).Value = 1.183405100000001
$ws.Cells.Item(16,15).Value = 1.315509199999998
$ws.Cells.Item(16,16).Value = 0.8288598849999991
$ws.Cells.Item(16,17).Value = 1.090487835
$ws.Cells.Item(16,18).Value = 0.9470416233333331
$ws.Cells.Item(16,19).Value = 0.8410620799999996
$ws.Cells.Item(16,20).Value = 0.9470416233333331
$ws.Cells.Item(16,21).Value = 0.926647835
$ws.Cells.Item(16,22).Value = 0.9779992880000001
$ws.Cells.Item(16,23).Value = 0.9504921187500001
$ws.Cells.Item(17,3).Value = 0.4511150999999999
$ws.Cells.Item(17,4).Value = 1.3998199
$ws.Cells.Item(17,5).Value = 0.82918194
$ws.Cells.Item(17,6).Value = 1.2265952
$ws.Cells.Item(17,7).Value = 2.157425
$ws.Cells.Item(17,8).Value = 2.157425
$ws.Cells.Item(17,9).Value = 2.157425
$ws.Cells.Item(17,10).Value = 0.80455108
$ws.Cells.Item(17,11).Value = 0.23327132
$ws.Cells.Item(17,12).Value = 0.23327132
$ws.Cells.Item(17,13).Value = 0.62072248
$ws.Cells.Item(17,14).Value = 2.157425
$ws.Cells.Item(17,15).Value = 0.80455108
$ws.Cells.Item(17,16).Value = 0.5189112
$ws.Cells.Item(17,17).Value = 0.81686651
$ws.Cells.Item(17,18).Value = 1.065082466666667
$ws.Cells.Item(17,19).Value = 0.62233478
$ws.Cells.Item(17,20).Value = 1.065082466666667
$ws.Cells.Item(17,21).Value = 1.006107335
$ws.Cells.Item(17,22).Value = 1.236370868
$ws.Cells.Item(17,23).Value = 0.9653352525000001
$ws.Cells.Item(18,3).Value = 0.3768648800000001
$ws.Cells.Item(18,4).Value = 1.2071827
$ws.Cells.Item(18,5).Value = 0.72413805
$ws.Cells.Item(18,6).Value = 1.3825474
$ws.Cells.Item(18,7).Value = 1.7583183
$ws.Cells.Item(18,8).Value = 1.7583183
$ws.Cells.Item(18,9).Value = 1.7583183
$ws.Cells.Item(18,10).Value = 1.2625809
$ws.Cells.Item(18,11).Value = 0.16077021
$ws.Cells.Item(18,12).Value = 0.16077021
$ws.Cells.Item(18,13).Value = 0.7639211299999999
$ws.Cells.Item(18,14).Value = 1.7583183
$ws.Cells.Item(18,15).Value = 1.2625809
$ws.Cells.Item(18,16).Value = 0.7116755549999999
$ws.Cells.Item(18,17).Value = 0.9933594749999999
$ws.Cells.Item(18,18).Value = 1.06055647
$ws.Cells.Item(18,19).Value = 0.7158297199999999
$ws.Cells.Item(18,20).Value = 1.06055647
$ws.Cells.Item(18,21).Value = 0.9764518649999999
$ws.Cells.Item(18,22).Value = 1.132825152
$ws.Cells.Item(18,23).Value = 0.95454044625
$ws.Cells.Item(19,3).Value = 0.8204790199999999
$ws.Cells.Item(19,4).Value = 1.1092923
$ws.Cells.Item(19,5).Value = 0.98445046
$ws.Cells.Item(19,6).Value = 1.0739027
$ws.Cells.Item(19,7).Value = 1.18395
$ws.Cells.Item(19,8).Value = 1.18395
$ws.Cells.Item(19,9).Value = 1.18395
$ws.Cells.Item(19,10).Value = 1.007136
$ws.Cells.Item(19,11).Value = 0.6563165200000001
$ws.Cells.Item(19,12).Value = 0.6563165200000001
$ws.Cells.Item(19,13).Value = 0.9505221699999999
$ws.Cells.Item(19,14).Value = 1.18395
$ws.Cells.Item(19,15).Value = 1.007136
$ws.Cells.Item(19,16).Value = 0.8317262600000001
$ws.Cells.Item(19,17).Value = 0.9957932300000001
$ws.Cells.Item(19,18).Value = 0.9491341733333334
$ws.Cells.Item(19,19).Value = 0.8826343266666666
$ws.Cells.Item(19,20).Value = 0.9491341733333334
$ws.Cells.Item(19,21).Value = 0.9579632450000001
$ws.Cells.Item(19,22).Value = 1.003160596
$ws.Cells.Item(19,23).Value = 0.97325614625
$ws.Cells.Item(20,3).Value = 0.8671396079452055
$ws.Cells.Item(20,4).Value = 1.071806261643836
$ws.Cells.Item(20,5).Value = 0.9539852597260273
$ws.Cells.Item(20,6).Value = 1.064457720547946
$ws.Cells.Item(20,7).Value = 1.219908410410959
$ws.Cells.Item(20,8).Value = 1.219908410410959
$ws.Cells.Item(20,9).Value = 1.219908410410959
$ws.Cells.Item(20,10).Value = 1.000584686575342
$ws.Cells.Item(20,11).Value = 0.8068983608219181
$ws.Cells.Item(20,12).Value = 0.8068983608219181
$ws.Cells.Item(20,13).Value = 0.9334656342465755
$ws.Cells.Item(20,14).Value = 1.219908410410959
$ws.Cells.Item(20,15).Value = 1.000584686575342
$ws.Cells.Item(20,16).Value = 0.9037415236986301
$ws.Cells.Item(20,17).Value = 0.9772849731506847
$ws.Cells.Item(20,18).Value = 1.009130485936073
$ws.Cells.Item(20,19).Value = 0.9204894357077625
$ws.Cells.Item(20,20).Value = 1.009130485936073
$ws.Cells.Item(20,21).Value = 0.9953441793835618
$ws.Cells.Item(20,22).Value = 1.040257025589041
$ws.Cells.Item(20,23).Value = 0.9897807427397263
$ws.Cells.Item(21,3).Value = 0.7051555784210526
$ws.Cells.Item(21,4).Value = 1.060169268947368
$ws.Cells.Item(21,5).Value = 0.9047324863157896
$ws.Cells.Item(21,6).Value = 1.190464578947368
$ws.Cells.Item(21,7).Value = 1.12050397
$ws.Cells.Item(21,8).Value = 1.12050397
$ws.Cells.Item(21,9).Value = 1.12050397
$ws.Cells.Item(21,10).Value = 1.232643878947368
$ws.Cells.Item(21,11).Value = 0.5157660921052631
$ws.Cells.Item(21,12).Value = 0.5157660921052631
$ws.Cells.Item(21,13).Value = 0.974875287894737
$ws.Cells.Item(21,14).Value = 1.12050397
$ws.Cells.Item(21,15).Value = 1.232643878947368
$ws.Cells.Item(21,16).Value = 0.8742049855263156
$ws.Cells.Item(21,17).Value = 1.068688182631579
$ws.Cells.Item(21,18).Value = 0.9563046470175438
$ws.Cells.Item(21,19).Value = 0.884380819122807
$ws.Cells.Item(21,20).Value = 0.9563046470175438
$ws.Cells.Item(21,21).Value = 0.9434116068421052
$ws.Cells.Item(21,22).Value = 0.9788300794736842
$ws.Cells.Item(21,23).Value = 0.9630388926973684
$ws.Cells.Item(22,3).Value = 0.7446285189473685
$ws.Cells.Item(22,4).Value = 1.084647864210526
$ws.Cells.Item(22,5).Value = 0.8995358031578948
$ws.Cells.Item(22,6).Value = 1.158401121578948
$ws.Cells.Item(22,7).Value = 1.256881638947368
$ws.Cells.Item(22,8).Value = 1.256881638947368
$ws.Cells.Item(22,9).Value = 1.256881638947368
$ws.Cells.Item(22,10).Value = 1.122981076842105
$ws.Cells.Item(22,11).Value = 0.637503872631579
$ws.Cells.Item(22,12).Value = 0.637503872631579
$ws.Cells.Item(22,13).Value = 0.9163306931578947
$ws.Cells.Item(22,14).Value = 1.256881638947368
$ws.Cells.Item(22,15).Value = 1.122981076842105
$ws.Cells.Item(22,16).Value = 0.8802424747368423
$ws.Cells.Item(22,17).Value = 1.01125844
$ws.Cells.Item(22,18).Value = 1.005788862807018
$ws.Cells.Item(22,19).Value = 0.8866735842105266
$ws.Cells.Item(22,20).Value = 1.005788862807018
$ws.Cells.Item(22,21).Value = 0.9792255978947368
$ws.Cells.Item(22,22).Value = 1.034756806105263
$ws.Cells.Item(22,23).Value = 0.9776138236842107
$ws.Cells.Item(23,3).Value = 0.9443093510656936
$ws.Cells.Item(23,4).Value = 0.9400114006175166
$ws.Cells.Item(23,5).Value = 0.9839295210114044
$ws.Cells.Item(23,6).Value = 1.064927450800309
$ws.Cells.Item(23,7).Value = 0.7861127358559966
$ws.Cells.Item(23,8).Value = 0.7861127358559966
$ws.Cells.Item(23,9).Value = 0.7861127358559966
$ws.Cells.Item(23,10).Value = 1.199091074229437
$ws.Cells.Item(23,11).Value = 0.8654670683000033
$ws.Cells.Item(23,12).Value = 0.8654670683000033
$ws.Cells.Item(23,13).Value = 1.083132371393852
$ws.Cells.Item(23,14).Value = 0.7861127358559966
$ws.Cells.Item(23,15).Value = 1.199091074229437
$ws.Cells.Item(23,16).Value = 1.03227907126472
$ws.Cells.Item(23,17).Value = 1.091510297620421
$ws.Cells.Item(23,18).Value = 0.950223626128479
$ws.Cells.Item(23,19).Value = 1.016162554513615
$ws.Cells.Item(23,20).Value = 0.950223626128479
$ws.Cells.Item(23,21).Value = 0.9586500998492103
$ws.Cells.Item(23,22).Value = 0.9241426270505675
$ws.Cells.Item(23,23).Value = 0.9833726216592764
$ws.Cells.Item(24,3).Value = 1.019218186962218
$ws.Cells.Item(24,4).Value = 1.025432234718397
$ws.Cells.Item(24,5).Value = 1.003998577029389
$ws.Cells.Item(24,6).Value = 0.9768914029675525
$ws.Cells.Item(24,7).Value = 1.090910876987596
$ws.Cells.Item(24,8).Value = 1.090910876987596
$ws.Cells.Item(24,9).Value = 1.090910876987596
$ws.Cells.Item(24,10).Value = 0.9213753074936379
$ws.Cells.Item(24,11).Value = 1.053476317894017
$ws.Cells.Item(24,12).Value = 1.053476317894017
$ws.Cells.Item(24,13).Value = 0.9633449958563002
$ws.Cells.Item(24,14).Value = 1.090910876987596
$ws.Cells.Item(24,15).Value = 0.9213753074936379
$ws.Cells.Item(24,16).Value = 0.9874258126938273
$ws.Cells.Item(24,17).Value = 0.9626869422615136
$ws.Cells.Item(24,18).Value = 1.021920834125084
$ws.Cells.Item(24,19).Value = 0.992950067472348
$ws.Cells.Item(24,20).Value = 1.021920834125084
$ws.Cells.Item(24,21).Value = 1.01744026985116
$ws.Cells.Item(24,22).Value = 1.032134391278448
$ws.Cells.Item(24,23).Value = 1.006830987488639
$ws.Cells.Item(25,3).Value = 1.121099826466081
$ws.Cells.Item(25,4).Value = 1.002430747304343
$ws.Cells.Item(25,5).Value = 1.041978687961073
$ws.Cells.Item(25,6).Value = 0.910020293789749
$ws.Cells.Item(25,7).Value = 1.029052815545128
$ws.Cells.Item(25,8).Value = 1.029052815545128
$ws.Cells.Item(25,9).Value = 1.029052815545128
$ws.Cells.Item(25,10).Value = 0.8480313563145121
$ws.Cells.Item(25,11).Value = 1.208438615235842
$ws.Cells.Item(25,12).Value = 1.208438615235842
$ws.Cells.Item(25,13).Value = 0.9811293604538415
$ws.Cells.Item(25,14).Value = 1.029052815545128
$ws.Cells.Item(25,15).Value = 0.8480313563145121
$ws.Cells.Item(25,16).Value = 1.028234985775177
$ws.Cells.Item(25,17).Value = 0.9450050221377926
$ws.Cells.Item(25,18).Value = 1.028507595698494
$ws.Cells.Item(25,19).Value = 1.032816219837143
$ws.Cells.Item(25,20).Value = 1.028507595698494
$ws.Cells.Item(25,21).Value = 1.031875368764139
$ws.Cells.Item(25,22).Value = 1.031310858120337
$ws.Cells.Item(25,23).Value = 1.017772712883821
$ws.Cells.Item(26,3).Value = 1.011653872214887
$ws.Cells.Item(26,4).Value = 0.9731921456920544
$ws.Cells.Item(26,5).Value = 1.00103054220385
$ws.Cells.Item(26,6).Value = 1.002846928100491
$ws.Cells.Item(26,7).Value = 0.9253501863210959
$ws.Cells.Item(26,8).Value = 0.9253501863210959
$ws.Cells.Item(26,9).Value = 0.9253501863210959
$ws.Cells.Item(26,10).Value = 1.041642391419102
$ws.Cells.Item(26,11).Value = 1.011123206271381
$ws.Cells.Item(26,12).Value = 1.011123206271381
$ws.Cells.Item(26,13).Value = 1.027129672452167
$ws.Cells.Item(26,14).Value = 0.9253501863210959
$ws.Cells.Item(26,15).Value = 1.041642391419102
$ws.Cells.Item(26,16).Value = 1.026382798845241
$ws.Cells.Item(26,17).Value = 1.021336466811476
$ws.Cells.Item(26,18).Value = 0.9927052613371928
$ws.Cells.Item(26,19).Value = 1.017932046631444
$ws.Cells.Item(26,20).Value = 0.9927052613371928
$ws.Cells.Item(26,21).Value = 0.994786581553857
$ws.Cells.Item(26,22).Value = 0.9808993025073048
$ws.Cells.Item(26,23).Value = 0.9992461180843784
$ws.Cells.Item(27,3).Value = 1.048559597174507
$ws.Cells.Item(27,4).Value = 1.029425751830572
$ws.Cells.Item(27,5).Value = 1.009964475457077
$ws.Cells.Item(27,6).Value = 0.95488794712036
$ws.Cells.Item(27,7).Value = 1.120695045832896
$ws.Cells.Item(27,8).Value = 1.120695045832896
$ws.Cells.Item(27,9).Value = 1.120695045832896
$ws.Cells.Item(27,10).Value = 0.8749414482001507
$ws.Cells.Item(27,11).Value = 1.114912747444406
$ws.Cells.Item(27,12).Value = 1.114912747444406
$ws.Cells.Item(27,13).Value = 0.9495843473075586
$ws.Cells.Item(27,14).Value = 1.120695045832896
$ws.Cells.Item(27,15).Value = 0.8749414482001507
$ws.Cells.Item(27,16).Value = 0.9949270978222784
$ws.Cells.Item(27,17).Value = 0.9424529618286137
$ws.Cells.Item(27,18).Value = 1.036849747159151
$ws.Cells.Item(27,19).Value = 0.9999395570338777
$ws.Cells.Item(27,20).Value = 1.036849747159151
$ws.Cells.Item(27,21).Value = 1.030128429233632
$ws.Cells.Item(27,22).Value = 1.048241752553485
$ws.Cells.Item(27,23).Value = 1.012871420045941
$ws.Cells.Item(28,3).Value = 0.9915000196836785
$ws.Cells.Item(28,4).Value = 0.9849426755539294
$ws.Cells.Item(28,5).Value = 1.003680113288748
$ws.Cells.Item(28,6).Value = 1.007665716669143
$ws.Cells.Item(28,7).Value = 0.940333488648821
$ws.Cells.Item(28,8).Value = 0.940333488648821
$ws.Cells.Item(28,9).Value = 0.940333488648821
$ws.Cells.Item(28,10).Value = 1.042119200074225
$ws.Cells.Item(28,11).Value = 0.9587770631887965
$ws.Cells.Item(28,12).Value = 0.9587770631887965
$ws.Cells.Item(28,13).Value = 1.028784504739448
$ws.Cells.Item(28,14).Value = 0.940333488648821
$ws.Cells.Item(28,15).Value = 1.042119200074225
$ws.Cells.Item(28,16).Value = 1.000448131631511
$ws.Cells.Item(28,17).Value = 1.022899656681487
$ws.Cells.Item(28,18).Value = 0.9804099173039477
$ws.Cells.Item(28,19).Value = 1.00152545885059
$ws.Cells.Item(28,20).Value = 0.9804099173039477
$ws.Cells.Item(28,21).Value = 0.9862274663001477
$ws.Cells.Item(28,22).Value = 0.9770486707698824
$ws.Cells.Item(28,23).Value = 0.9947253477308489
$ws.Cells.Item(29,3).Value = 0.9784866171187694
$ws.Cells.Item(29,4).Value = 1.008370682216592
$ws.Cells.Item(29,5).Value = 0.988517652051378
$ws.Cells.Item(29,6).Value = 1.016219640118079
$ws.Cells.Item(29,7).Value = 1.025290401903346
$ws.Cells.Item(29,8).Value = 1.025290401903346
$ws.Cells.Item(29,9).Value = 1.025290401903346
$ws.Cells.Item(29,10).Value = 1.010225248292476
$ws.Cells.Item(29,11).Value = 0.9810484054479265
$ws.Cells.Item(29,12).Value = 0.9810484054479265
$ws.Cells.Item(29,13).Value = 0.9864832632323143
$ws.Cells.Item(29,14).Value = 1.025290401903346
$ws.Cells.Item(29,15).Value = 1.010225248292476
$ws.Cells.Item(29,16).Value = 0.9956368268702014
$ws.Cells.Item(29,17).Value = 0.9993714501719271
$ws.Cells.Item(29,18).Value = 1.00552135188125
$ws.Cells.Item(29,19).Value = 0.9932637685972603
$ws.Cells.Item(29,20).Value = 1.00552135188125
$ws.Cells.Item(29,21).Value = 1.001270426923782
$ws.Cells.Item(29,22).Value = 1.006074421919695
$ws.Cells.Item(29,23).Value = 0.9993302387976102
$ws.Cells.Item(30,3).Value = 0.9735090803331943
$ws.Cells.Item(30,4).Value = 1.006673564291857
$ws.Cells.Item(30,5).Value = 0.9844439554057631
$ws.Cells.Item(30,6).Value = 1.017228998004731
$ws.Cells.Item(30,7).Value = 1.0388049669224
$ws.Cells.Item(30,8).Value = 1.0388049669224
$ws.Cells.Item(30,9).Value = 1.0388049669224
$ws.Cells.Item(30,10).Value = 1.011750884460101
$ws.Cells.Item(30,11).Value = 0.9698449742054086
$ws.Cells.Item(30,12).Value = 0.9698449742054086
$ws.Cells.Item(30,13).Value = 0.9886292071114462
$ws.Cells.Item(30,14).Value = 1.0388049669224
$ws.Cells.Item(30,15).Value = 1.011750884460101
$ws.Cells.Item(30,16).Value = 0.9907979293327547
$ws.Cells.Item(30,17).Value = 0.998097419932932
$ws.Cells.Item(30,18).Value = 1.006800275195969
$ws.Cells.Item(30,19).Value = 0.9886799380237575
$ws.Cells.Item(30,20).Value = 1.00680027519597
$ws.Cells.Item(30,21).Value = 1.001211195248418
$ws.Cells.Item(30,22).Value = 1.008729949583214
$ws.Cells.Item(30,23).Value = 0.9988607038418625
$ws.Cells.Item(31,3).Value = 1.052694225096709
$ws.Cells.Item(31,4).Value = 1.011816659594402
$ws.Cells.Item(31,5).Value = 1.035770685897295
$ws.Cells.Item(31,6).Value = 0.9509163744030563
$ws.Cells.Item(31,7).Value = 0.9991287781189778
$ws.Cells.Item(31,8).Value = 0.9991287781189778
$ws.Cells.Item(31,9).Value = 0.9991287781189778
$ws.Cells.Item(31,10).Value = 0.9181375634241195
$ws.Cells.Item(31,11).Value = 1.057919932967258
$ws.Cells.Item(31,12).Value = 1.057919932967258
$ws.Cells.Item(31,13).Value = 1.000980766683897
$ws.Cells.Item(31,14).Value = 0.9991287781189778
$ws.Cells.Item(31,15).Value = 0.9181375634241195
$ws.Cells.Item(31,16).Value = 0.9880287481956886
$ws.Cells.Item(31,17).Value = 0.9769541246607072
$ws.Cells.Item(31,18).Value = 0.9917287581701184
$ws.Cells.Item(31,19).Value = 1.003942727429558
$ws.Cells.Item(31,20).Value = 0.9917287581701184
$ws.Cells.Item(31,21).Value = 1.002739240101913
$ws.Cells.Item(31,22).Value = 1.002017147705326
$ws.Cells.Item(31,23).Value = 1.003420623273214

# --- Step 7: copy the column-A formatting (bold, border, centered) onto the two new rows
$ws.Range("A29").Copy() | Out-Null
$ws.Range("A30:A31").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
